$d = $word.ActiveDocument

$d.Content.Find.Execute("117÷2=58, 1", $true, $false, $false, $false, $false, $true, 1, $false, "487÷7=69, 4", 2) | Out-Null
$d.Content.Find.Execute("901÷8=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "488÷2=244, 0", 2) | Out-Null
$d.Content.Find.Execute("184÷3=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "358÷9=39, 7", 2) | Out-Null
$d.Content.Find.Execute("934÷4=233, 2", $true, $false, $false, $false, $false, $true, 1, $false, "495÷6=82, 3", 2) | Out-Null
$d.Content.Find.Execute("479÷7=68, 3", $true, $false, $false, $false, $false, $true, 1, $false, "414÷8=51, 6", 2) | Out-Null
$d.Content.Find.Execute("493÷7=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "615÷9=68, 3", 2) | Out-Null
$d.Content.Find.Execute("963÷6=160, 3", $true, $false, $false, $false, $false, $true, 1, $false, "749÷4=187, 1", 2) | Out-Null
$d.Content.Find.Execute("110÷3=36, 2", $true, $false, $false, $false, $false, $true, 1, $false, "390÷9=43, 3", 2) | Out-Null
$d.Content.Find.Execute("463÷2=231, 1", $true, $false, $false, $false, $false, $true, 1, $false, "822÷9=91, 3", 2) | Out-Null
$d.Content.Find.Execute("754÷5=150, 4", $true, $false, $false, $false, $false, $true, 1, $false, "786÷8=98, 2", 2) | Out-Null
$d.Content.Find.Execute("770÷9=85, 5", $true, $false, $false, $false, $false, $true, 1, $false, "393÷6=65, 3", 2) | Out-Null
$d.Content.Find.Execute("669÷2=334, 1", $true, $false, $false, $false, $false, $true, 1, $false, "253÷2=126, 1", 2) | Out-Null
$d.Content.Find.Execute("572÷8=71, 4", $true, $false, $false, $false, $false, $true, 1, $false, "774÷3=258, 0", 2) | Out-Null
$d.Content.Find.Execute("991÷3=330, 1", $true, $false, $false, $false, $false, $true, 1, $false, "285÷6=47, 3", 2) | Out-Null
$d.Content.Find.Execute("687÷2=343, 1", $true, $false, $false, $false, $false, $true, 1, $false, "426÷6=71, 0", 2) | Out-Null
$d.Content.Find.Execute("802÷6=133, 4", $true, $false, $false, $false, $false, $true, 1, $false, "360÷6=60, 0", 2) | Out-Null
$d.Content.Find.Execute("489÷2=244, 1", $true, $false, $false, $false, $false, $true, 1, $false, "623÷6=103, 5", 2) | Out-Null
$d.Content.Find.Execute("277÷9=30, 7", $true, $false, $false, $false, $false, $true, 1, $false, "284÷3=94, 2", 2) | Out-Null
$d.Content.Find.Execute("711÷2=355, 1", $true, $false, $false, $false, $false, $true, 1, $false, "112÷3=37, 1", 2) | Out-Null
$d.Content.Find.Execute("634÷9=70, 4", $true, $false, $false, $false, $false, $true, 1, $false, "120÷7=17, 1", 2) | Out-Null
$d.Content.Find.Execute("338÷3=112, 2", $true, $false, $false, $false, $false, $true, 1, $false, "399÷5=79, 4", 2) | Out-Null
$d.Content.Find.Execute("938÷4=234, 2", $true, $false, $false, $false, $false, $true, 1, $false, "613÷5=122, 3", 2) | Out-Null
$d.Content.Find.Execute("536÷8=67, 0", $true, $false, $false, $false, $false, $true, 1, $false, "561÷3=187, 0", 2) | Out-Null
$d.Content.Find.Execute("346÷4=86, 2", $true, $false, $false, $false, $false, $true, 1, $false, "800÷7=114, 2", 2) | Out-Null
$d.Content.Find.Execute("234÷2=117, 0", $true, $false, $false, $false, $false, $true, 1, $false, "862÷2=431, 0", 2) | Out-Null
